$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to snake_case identifiers
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Convert state / municipality / "TOTAL" text from upper case to title case
$ws.Range("A2").Value = "Baja California"
$ws.Range("B2").Value = "Mexicali"
$ws.Range("B3").Value = "Tecate"
$ws.Range("B4").Value = "Tijuana"
$ws.Range("B5").Value = "Total"
$ws.Range("A6").Value = "Campeche"
$ws.Range("B6").Value = "Champotón"
$ws.Range("B7").Value = "Total"
$ws.Range("A8").Value = "Chiapas"
$ws.Range("B8").Value = "Acapetahua"
$ws.Range("B9").Value = "Jiquipilas"
$ws.Range("B10").Value = "Mapastepec"
$ws.Range("B11").Value = "Pijijiapan"
$ws.Range("B12").Value = "Total"
$ws.Range("A13").Value = "Chihuahua"
$ws.Range("B13").Value = "Chihuahua"
$ws.Range("B14").Value = "Cuauhtémoc"
$ws.Range("B15").Value = "Total"
$ws.Range("A16").Value = "Ciudad De México"
$ws.Range("B16").Value = "Azcapotzalco"
$ws.Range("B17").Value = "Benito Juárez"
$ws.Range("B18").Value = "Coyoacán"
$ws.Range("B19").Value = "Cuauhtémoc"
$ws.Range("B20").Value = "Gustavo A. Madero"
$ws.Range("B21").Value = "Iztacalco"
$ws.Range("B22").Value = "Iztapalapa"
$ws.Range("B23").Value = "Miguel Hidalgo"
$ws.Range("B24").Value = "No Se Registró El Municipio/Condado/Alcaldía De Nacimiento"
$ws.Range("B25").Value = "Venustiano Carranza"
$ws.Range("B26").Value = "Álvaro Obregón"
$ws.Range("B27").Value = "Total"
$ws.Range("A28").Value = "Durango"
$ws.Range("B28").Value = "Nazas"
$ws.Range("B29").Value = "Total"
$ws.Range("A30").Value = "Estado De México"
$ws.Range("B30").Value = "Chicoloapan"
$ws.Range("B31").Value = "Nezahualcóyotl"
$ws.Range("B32").Value = "Ozumba"
$ws.Range("B33").Value = "Sultepec"
$ws.Range("B34").Value = "Tenancingo"
$ws.Range("B35").Value = "Toluca"
$ws.Range("B36").Value = "Valle De Chalco Solidaridad"
$ws.Range("B37").Value = "Total"
$ws.Range("A38").Value = "Guanajuato"
$ws.Range("B38").Value = "Salvatierra"
$ws.Range("B39").Value = "Total"
$ws.Range("A40").Value = "Guerrero"
$ws.Range("B40").Value = "Buenavista De Cuéllar"
$ws.Range("B41").Value = "Pungarabato"
$ws.Range("B42").Value = "Total"
$ws.Range("A43").Value = "Hidalgo"
$ws.Range("B43").Value = "Acaxochitlán"
$ws.Range("B44").Value = "Actopan"
$ws.Range("B45").Value = "Atotonilco El Grande"
$ws.Range("B46").Value = "Metepec"
$ws.Range("B47").Value = "Tulancingo De Bravo"
$ws.Range("B48").Value = "Total"
$ws.Range("A49").Value = "Jalisco"
$ws.Range("B49").Value = "Arandas"
$ws.Range("B50").Value = "Autlán De Navarro"
$ws.Range("B51").Value = "Casimiro Castillo"
$ws.Range("B52").Value = "Guadalajara"
$ws.Range("B53").Value = "La Huerta"
$ws.Range("B54").Value = "Magdalena"
$ws.Range("B55").Value = "San Pedro Tlaquepaque"
$ws.Range("B56").Value = "Tonalá"
$ws.Range("B57").Value = "Tuxpan"
$ws.Range("B58").Value = "Total"
$ws.Range("A59").Value = "Michoacán De Ocampo"
$ws.Range("B59").Value = "La Piedad"
$ws.Range("B60").Value = "Los Reyes"
$ws.Range("B61").Value = "Puruándiro"
$ws.Range("B62").Value = "Salvador Escalante"
$ws.Range("B63").Value = "Villamar"
$ws.Range("B64").Value = "Zacapu"
$ws.Range("B65").Value = "Zamora"
$ws.Range("B66").Value = "Total"
$ws.Range("A67").Value = "Nayarit"
$ws.Range("B67").Value = "La Yesca"
$ws.Range("B68").Value = "San Blas"
$ws.Range("B69").Value = "Santa María Del Oro"
$ws.Range("B70").Value = "Tepic"
$ws.Range("B71").Value = "Total"
$ws.Range("A72").Value = "Oaxaca"
$ws.Range("B72").Value = "Chalcatongo De Hidalgo"
$ws.Range("B73").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B74").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B75").Value = "Putla Villa De Guerrero"
$ws.Range("B76").Value = "San Andrés Cabecera Nueva"
$ws.Range("B77").Value = "San Francisco Del Mar"
$ws.Range("B78").Value = "Santa María Cortijo"
$ws.Range("B79").Value = "Santiago Jamiltepec"
$ws.Range("B80").Value = "Santo Domingo Tehuantepec"
$ws.Range("B81").Value = "Santo Tomás Ocotepec"
$ws.Range("B82").Value = "Total"
$ws.Range("A83").Value = "Puebla"
$ws.Range("B83").Value = "Libres"
$ws.Range("B84").Value = "Los Reyes De Juárez"
$ws.Range("B85").Value = "Tetela De Ocampo"
$ws.Range("B86").Value = "Tlachichuca"
$ws.Range("B87").Value = "Total"
$ws.Range("A88").Value = "Querétaro"
$ws.Range("B88").Value = "Amealco De Bonfil"
$ws.Range("B89").Value = "Cadereyta De Montes"
$ws.Range("B90").Value = "Total"
$ws.Range("A91").Value = "San Luis Potosí"
$ws.Range("B91").Value = "Tamazunchale"
$ws.Range("B92").Value = "Total"
$ws.Range("A93").Value = "Sinaloa"
$ws.Range("B93").Value = "El Fuerte"
$ws.Range("B94").Value = "Total"
$ws.Range("A95").Value = "Sonora"
$ws.Range("B95").Value = "Hermosillo"
$ws.Range("B96").Value = "Total"
$ws.Range("A97").Value = "Tabasco"
$ws.Range("B97").Value = "Jalapa"
$ws.Range("B98").Value = "Total"
$ws.Range("A99").Value = "Tlaxcala"
$ws.Range("B99").Value = "Huamantla"
$ws.Range("B100").Value = "Total"
$ws.Range("A101").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B101").Value = "Papantla"
$ws.Range("B102").Value = "Total"
$ws.Range("A103").Value = "Zacatecas"
$ws.Range("B103").Value = "Sain Alto"
$ws.Range("B104").Value = "Valparaíso"
$ws.Range("B105").Value = "Total"
$ws.Range("A106").Value = "Total"

# Remove the trailing metadata/footer rows (108-112), shrinking the used range
# from A1:D112 down to A1:D106
$ws.Range("A108:D112").EntireRow.Delete()
